$wb = $excel.ActiveWorkbook

# Trade #7 closed at 2026-02-16 22:57:03 - base_strategy DOWN +0.000%
# Append the new trade row to both the "All Trades" sheet and the
# strategy-specific "base_strategy" sheet.

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 8

    # Leading apostrophe forces text interpretation (quote-prefix) so
    # Excel doesn't auto-convert the ISO date string to a date serial,
    # and so the otherwise-empty cells are still materialized (an
    # assignment of a bare "" does not create a cell at all).
    $ws.Cells.Item($row, 1).Value = 7
    $ws.Cells.Item($row, 2).Value = "'2026-02-16"
    $ws.Cells.Item($row, 3).Value = "22:57:03"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.5
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 17).Value = 0
}
